$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'98.854.37"
$ws.Range("E2").Value = "'  +0.89%  "
$ws.Range("D3").Value = "'3.342.57"
$ws.Range("E3").Value = "'  +6.33%  "
$ws.Range("E4").Value = "'  -0.06%  "
$ws.Range("D5").Value = "'258.61"
$ws.Range("E5").Value = "'  +7.31%  "
$ws.Range("D6").Value = "'625.09"
$ws.Range("E6").Value = "'  +2.55%  "
$ws.Range("D7").Value = "'1.44"
$ws.Range("E7").Value = "'  +29.82%  "
$ws.Range("D8").Value = "'0.393"
$ws.Range("E8").Value = "'  +2.84%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "'  -0.07%  "
$ws.Range("D10").Value = "'0.889"
$ws.Range("E10").Value = "'  +12.90%  "
$ws.Range("D11").Value = "'3.339.80"
$ws.Range("E12").Value = "'  +0.19%  "
$ws.Range("D13").Value = "'37.55"
$ws.Range("E13").Value = "'  +10.69%  "
$ws.Range("D14").Value = "'98.532.37"
$ws.Range("E14").Value = "'  +0.99%  "
$ws.Range("E15").Value = "'  +4.16%  "
$ws.Range("D16").Value = "'3.961.82"
$ws.Range("E16").Value = "'  +6.31%  "
$ws.Range("D17").Value = "'5.53"
$ws.Range("E17").Value = "'  +2.14%  "
$ws.Range("D18").Value = "'3.340.92"
$ws.Range("E18").Value = "'  +6.35%  "
$ws.Range("E19").Value = "'  +3.78%  "
$ws.Range("D20").Value = "'15.26"
$ws.Range("E20").Value = "'  +4.88%  "
$ws.Range("D21").Value = "'491.06"
$ws.Range("E21").Value = "'  -5.90%  "
$ws.Range("D22").Value = "'6.12"
$ws.Range("E22").Value = "'  +7.07%  "
$ws.Range("E23").Value = "'  +10.17%  "
$ws.Range("D24").Value = "'9.40"
$ws.Range("E24").Value = "'  +6.97%  "
$ws.Range("D25").Value = "'5.64"
$ws.Range("E25").Value = "'  +3.35%  "
$ws.Range("D26").Value = "'89.47"
$ws.Range("E26").Value = "'  +1.07%  "
$ws.Range("D27").Value = "'11.96"
$ws.Range("E27").Value = "'  +3.01%  "
$ws.Range("E28").Value = "'  +6.13%  "
$ws.Range("D29").Value = "'0.296"
$ws.Range("E29").Value = "'  +24.48%  "
$ws.Range("E30").Value = "'  -0.14%  "
$ws.Range("D31").Value = "'0.193"
$ws.Range("E31").Value = "'  +9.36%  "
$ws.Range("D32").Value = "'0.138"
$ws.Range("E32").Value = "'  +12.86%  "
$ws.Range("D33").Value = "'9.76"
$ws.Range("E33").Value = "'  +9.20%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "'  +0.31%  "
$ws.Range("D35").Value = "'28.23"
$ws.Range("E35").Value = "'  +6.08%  "
$ws.Range("B36").Value = "'RenderToken"
$ws.Range("C36").Value = "'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D36").Value = "'7.29"
$ws.Range("E36").Value = "'  +1.15%  "
$ws.Range("B37").Value = "'Kaspa"
$ws.Range("C37").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.150"
$ws.Range("E37").Value = "'  -0.95%  "
$ws.Range("D38").Value = "'1.95"
$ws.Range("E38").Value = "'  +4.13%  "
$ws.Range("D39").Value = "'501.15"
$ws.Range("E39").Value = "'  +7.88%  "
$ws.Range("E40").Value = "'  +5.96%  "
$ws.Range("D41").Value = "'24.90"
$ws.Range("E41").Value = "'  +2.29%  "
$ws.Range("E42").Value = "'  +3.98%  "
$ws.Range("D43").Value = "'3.66"
$ws.Range("E43").Value = "'  +4.12%  "
$ws.Range("D44").Value = "'3.34"
$ws.Range("E44").Value = "'  +7.48%  "
$ws.Range("E46").Value = "'  +12.11%  "
$ws.Range("D47").Value = "'159.67"
$ws.Range("E47").Value = "'  -1.52%  "
$ws.Range("D48").Value = "'1.95"
$ws.Range("E48").Value = "'  +1.18%  "
$ws.Range("D49").Value = "'0.850"
$ws.Range("E49").Value = "'  +8.21%  "
$ws.Range("D50").Value = "'4.65"
$ws.Range("E50").Value = "'  +3.11%  "
$ws.Range("E51").Value = "'  +4.06%  "
